# EPP Variable Installments T1 scenarios
#
# Applies the changes observed in the target diff:
#  - Summary sheet: update the current selection
#  - Repayment schedule sheet: drop the now-unused "O" column cell in
#    rows 3-15 (and the trailing empty cell in row 2), update selection
#  - Transactions sheet: new transaction IDs, narrower ID column, update
#    selection (this sheet stays the active/selected tab)

$wb = $excel.ActiveWorkbook

# --- Summary ---------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("F5").Select()

# --- Repayment schedule -----------------------------------------------
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")

# Row 2's trailing (empty) cell is dropped.
$wsSchedule.Range("P2").Clear()

# Rows 3-15 each lose their "O" column cell (column N and P stay put).
$wsSchedule.Range("O3:O15").Clear()

$wsSchedule.Range("G15").Select()

# --- Transactions -------------------------------------------------------
$wsTransactions = $wb.Worksheets.Item("Transactions")

$wsTransactions.Range("A2").Value = 1671
$wsTransactions.Range("A3").Value = 640
$wsTransactions.Range("A4").Value = 632

# Narrow the ID column from 4 to 3 characters.
$wsTransactions.Columns.Item(1).ColumnWidth = 2.1666666666666665

# Leave this as the active sheet/selection, matching tabSelected="1".
$wsTransactions.Range("J3").Select()
